# Regenerate the "K" (strikeouts) column (G) for the save_data sheet.
# The G column header is "K" but the previously-saved values were derived
# from the wrong source field (Strike# instead of real strikeout totals).
# This recalculates / rewrites the correct strikeout values ("s_vals")
# for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values, in row order, for rows 2..77 (data starts at row 2).
$sVals = @(
    1,1,0,1,0,0,1,1,1,0,
    1,0,2,3,2,1,2,1,0,1,
    3,1,2,1,1,5,3,3,2,1,
    3,2,2,3,2,1,1,1,2,2,
    1,2,3,1,3,4,4,2,1,1,
    2,2,4,3,2,1,1,0,0,2,
    1,1,2,1,3,1,2,2,1,5,
    2,1,3,0,3,2
)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
